$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "291.15") are not auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '22.390.18'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.570.74'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '291.15'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').Value = '0.3760'
$ws.Range('E7').Value = '  +1.81%  '
$ws.Range('D8').Value = '49.83'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = '0.3418'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('D10').Value = '0.07630'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').Value = '1.147'
$ws.Range('E11').Value = '  -2.15%  '
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '21.16'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('D14').Value = '6.001'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').Value = '6.931'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '1.571.97'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').Value = '90.31'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').Value = '0.06732'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').Value = '16.70'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('D22').Value = '6.210'
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').Value = '12.01'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '22.404.02'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '2.396'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').Value = '2.650'
$ws.Range('E26').Value = '  -12.47%  '
$ws.Range('D27').Value = '20.15'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').Value = '147.25'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('D30').Value = '126.68'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').Value = '1.745.13'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('D33').Value = '6.123'
$ws.Range('E33').Value = '  -2.88%  '
$ws.Range('D34').Value = '0.9828'
$ws.Range('E34').Value = '  -6.12%  '
$ws.Range('D35').Value = '10.06'
$ws.Range('E35').Value = '  -3.33%  '
$ws.Range('D36').Value = '0.08514'
$ws.Range('E36').Value = '  +0.65%  '
$ws.Range('D37').Value = '0.02535'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '1.359'
$ws.Range('E38').Value = '  +8.93%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('D40').Value = '0.06528'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('D41').Value = '5.411'
$ws.Range('E41').Value = '  -3.53%  '
$ws.Range('D42').Value = '11.42'
$ws.Range('E42').Value = '  -3.76%  '
$ws.Range('D43').Value = '0.6370'
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').Value = '14.05'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').Value = '3.801'
$ws.Range('D47').Value = '0.5962'
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('D48').Value = '1.289'
$ws.Range('E48').Value = '  +2.21%  '
$ws.Range('D49').Value = '2.092'
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('D50').Value = '124.43'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').Value = '0.07321'
$ws.Range('E51').Value = '  +0.25%  '
